$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily auction rows (10-Aug-2021 .. 02-Sep-2021) appended below the
# existing data (which ends at row 7).
#
# Columns: A=Serie(date), B=Cupo, C=Monto demandado, D=Total adjudicado,
#          E=Adjudicado bancos, F=Adjudicado AFP, G=Tasa interes
# Rows 9 and 18 only have Serie/Cupo/Total (no demand/rate was published),
# matching the sparse pattern already used by rows 5 and 6.
$rows = @(
    @{ Row=8;  Date="10-08-2021"; B=10000; C=22000; D=10000; E=10000; F=0;    G=2.1 },
    @{ Row=9;  Date="11-08-2021"; B=10000; D=0 },
    @{ Row=10; Date="12-08-2021"; B=10000; C=29000; D=10000; E=10000; F=0;    G=2.2 },
    @{ Row=11; Date="17-08-2021"; B=10000; C=29000; D=10000; E=10000; F=0;    G=2.21 },
    @{ Row=12; Date="18-08-2021"; B=10000; C=22000; D=10000; E=9000;  F=1000; G=2.22 },
    @{ Row=13; Date="19-08-2021"; B=10000; C=32000; D=10000; E=10000; F=0;    G=2 },
    @{ Row=14; Date="24-08-2021"; B=10000; C=37000; D=10000; E=10000; F=0;    G=2.15 },
    @{ Row=15; Date="25-08-2021"; B=10000; C=33000; D=10000; E=10000; F=0;    G=2.13 },
    @{ Row=16; Date="26-08-2021"; B=10000; C=22000; D=10000; E=9000;  F=1000; G=2.24 },
    @{ Row=17; Date="01-09-2021"; B=10000; C=19000; D=10000; E=10000; F=0;    G=3 },
    @{ Row=18; Date="02-09-2021"; B=10000; D=0 }
)

foreach ($entry in $rows) {
    $r = $entry.Row

    # Write the date as literal text (not an auto-converted date serial):
    # mark the cell as Text, assign the string, then drop the number-format
    # again so the cell keeps the workbook's default (unstyled) look.
    $aCell = $ws.Cells.Item($r, 1)
    $aCell.NumberFormat = "@"
    $aCell.Value = $entry.Date
    $aCell.ClearFormats()

    $ws.Cells.Item($r, 2).Value = $entry.B

    if ($entry.ContainsKey("C")) { $ws.Cells.Item($r, 3).Value = $entry.C }
    $ws.Cells.Item($r, 4).Value = $entry.D
    if ($entry.ContainsKey("E")) { $ws.Cells.Item($r, 5).Value = $entry.E }
    if ($entry.ContainsKey("F")) { $ws.Cells.Item($r, 6).Value = $entry.F }
    if ($entry.ContainsKey("G")) { $ws.Cells.Item($r, 7).Value = $entry.G }
}
